# Update "Comenzi Git" (Git commands) cheat-sheet worksheet:
# - add a new row (14) documenting "git branch" / "list of all branches"
# - move the highlighted-row formatting from the old last row (13) to the
#   new last row (14)
# - update the active selection to reflect where the user left off editing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 13 (currently the highlighted "last row") down into row 14 first,
# using Insert with CopyOrigin:=xlFormatFromLeftOrAbove (1) so the new row
# inherits row 13's cell formatting/style (the highlight fill) exactly.
$ws.Range("B13:C13").Copy()
$ws.Range("B14:C14").Insert(-4121, 1)

# Fill in the new row's content.
$ws.Range("B14").Value = "git branch"
$ws.Range("C14").Value = "list of all branches"

# The old last row (13) is no longer the last row, so it loses the
# highlight formatting, going back to the default/normal style.
$ws.Range("B13:C13").Style = "Normal"

# Reflect the cursor position left after the edit.
$ws.Range("C15").Select()
